# Rename sheets: sheet1 -> SPI3, Sheet2 -> SPI1
$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)

$ws1.Name = "SPI3"
$ws2.Name = "SPI1"

# The header text in column B of the second sheet ("SPI1") changes from
# "SPI1" to "SPI2".
$ws2.Range("B1").Value = "SPI2"

# Update the view/selection state:
#  - SPI1 (was the active/tabSelected sheet with topLeftCell A141 and a
#    selection of B152:B272) becomes inactive, scrolled back to the top,
#    with a simple single-cell selection at B1.
$ws2.Range("B1").Select()

#  - SPI3 (previously scrolled to A130 with a big selection) becomes the
#    active sheet, scrolled back to the top, with the selection/active
#    cell at I144.
$ws1.Activate()
$ws1.Range("I144").Select()
